$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    17 = 0.062853783369064331
    18 = 0.36311545968055725
    19 = 0.73902970552444458
    20 = 0.22195442020893097
    21 = 0.2039545476436615
    22 = 0.060609735548496246
    23 = 0.90907031297683716
    24 = 0.70409172773361206
    25 = 0.046301256865262985
    26 = 0.094824247062206268
    27 = 0.76530808210372925
    28 = 0.91357076168060303
    29 = 0.14032657444477081
    30 = 0.44705137610435486
    31 = 0.52541369199752808
    32 = 0.47234773635864258
    33 = 0.6911700963973999
    34 = 0.13295947015285492
    35 = 0.73065406084060669
    36 = 0.60111868381500244
    37 = 0.84555184841156006
    38 = 0.28247827291488647
    39 = 0.58727574348449707
    40 = 0.42868471145629883
    41 = 0.13384614884853363
    42 = 0.89993071556091309
    43 = 0.89632242918014526
    44 = 0.57144844532012939
    45 = 0.61684751510620117
    46 = 0.37054547667503357
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}
